$d = $word.ActiveDocument

# 1) Update the first three order lines (2023-04-09 entries -> new dates/amounts)
$d.Paragraphs.Item(3).Range.Text = "2023-05-20 - 11руб."
$d.Paragraphs.Item(4).Range.Text = "2023-05-12 - 6руб."
$d.Paragraphs.Item(5).Range.Text = "2023-05-01 - 4руб."

# 2) Remove the remaining order-line paragraphs for 2023-04-12 .. 2023-04-21
#    (paragraphs 6 through 15, i.e. everything between the three edited
#    order lines and the "Общая сумма" summary paragraph).
$startPara = $d.Paragraphs.Item(6)
$endPara = $d.Paragraphs.Item(15)
$r = $d.Range($startPara.Range.Start, $endPara.Range.End)
$r.Delete()

# 3) Update the total-sum summary line
$d.Paragraphs.Item(6).Range.Text = "Общая сумма заказов за прошлый месяц: 21руб."
